# Update the "Total fees accrued for the month" header (G1) to read
# "Fees accrued since last report" per the fixture refresh for FN-2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Fees accrued since last report"

# Move the active selection to G2, matching the saved cursor position in
# the updated fixture.
$ws.Range("G2").Select()
